$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1675.1111
$ws.Range("I96").Value = 1557
$ws.Range("J96").Value = 1769.6
$ws.Range("K96").Value = 4671
$ws.Range("L96").Value = 5308.799999999999
$ws.Range("M96").Value = -3298
$ws.Range("N96").Value = -8054.799999999999

$ws.Range("H111").Value = 4030
$ws.Range("I111").Value = 4030
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 12090
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -9023
$ws.Range("N111").ClearContents()

$ws.Range("H137").Value = 3826.6667
$ws.Range("I137").Value = 4113
$ws.Range("K137").Value = 12339
$ws.Range("M137").Value = -9789

$ws.Range("H138").Value = 3645.8286
$ws.Range("I138").Value = 3225.7273
$ws.Range("K138").Value = 9677.1819
$ws.Range("M138").Value = -4537.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1773.3077
$ws.Range("I122").Value = 1295.1111
$ws.Range("K122").Value = 3885.3333
$ws.Range("M122").Value = -1435.3333

$ws.Range("H130").Value = 29133.143
$ws.Range("J130").Value = 29133.143
$ws.Range("L130").Value = 29133.143
$ws.Range("N130").Value = -39173.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 940.6486
$ws.Range("I3").Value = 883.4722
$ws.Range("K3").Value = 883.4722
$ws.Range("M3").Value = -769.4722

$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 400
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 400
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = -260
$ws.Range("N8").Value = -880

$ws.Range("H11").Value = 1011.3333
$ws.Range("I11").Value = 743.1429000000001
$ws.Range("K11").Value = 743.1429000000001
$ws.Range("M11").Value = -603.1429000000001

$ws.Range("H12").Value = 6663
$ws.Range("J12").Value = 9925.75
$ws.Range("L12").Value = 9925.75
$ws.Range("N12").Value = -10261.75

$ws.Range("H105").Value = 3825.6
$ws.Range("I105").Value = 3784.3044
$ws.Range("J105").Value = 4300.5
$ws.Range("K105").Value = 3784.3044
$ws.Range("L105").Value = 4300.5
$ws.Range("M105").Value = -2037.3044
$ws.Range("N105").Value = -7794.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1768.5
$ws.Range("I22").Value = 1941.3334
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 1941.3334
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -1591.3334
$ws.Range("N22").Value = -1950

$ws.Range("H28").Value = 6081.1665
$ws.Range("J28").Value = 6081.1665
$ws.Range("L28").Value = 6081.1665
$ws.Range("N28").Value = -6571.1665

$ws.Range("H31").Value = 1994.6957
$ws.Range("I31").Value = 2195.353
$ws.Range("K31").Value = 2195.353
$ws.Range("M31").Value = -1900.353

$ws.Range("H34").Value = 1994.6957
$ws.Range("I34").Value = 2195.353
$ws.Range("K34").Value = 2195.353
$ws.Range("M34").Value = -1993.353

$ws.Range("H58").Value = 103222.3
$ws.Range("I58").Value = 145962.28
$ws.Range("K58").Value = 145962.28
$ws.Range("M58").Value = -145759.28

$ws.Range("H99").Value = 4181.5557
$ws.Range("J99").Value = 4557
$ws.Range("L99").Value = 4557
$ws.Range("N99").Value = -7553

$ws.Range("H116").Value = 29998
$ws.Range("J116").Value = 29998
$ws.Range("L116").Value = 29998
$ws.Range("N116").Value = -39176

$ws.Range("H126").Value = 4181.5557
$ws.Range("J126").Value = 4557
$ws.Range("L126").Value = 13671
$ws.Range("N126").Value = -18611

$ws.Range("H136").Value = 103222.3
$ws.Range("I136").Value = 145962.28
$ws.Range("K136").Value = 437886.84
$ws.Range("M136").Value = -435336.84

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1743.3636
$ws.Range("I7").Value = 294.5
$ws.Range("K7").Value = 883.5
$ws.Range("M7").Value = -771.5

$ws.Range("H13").Value = 863.3333
$ws.Range("I13").Value = 863.3333
$ws.Range("K13").Value = 2589.9999
$ws.Range("M13").Value = -2421.9999

$ws.Range("H34").Value = 1478.1428
$ws.Range("I34").Value = 1286.75
$ws.Range("J34").Value = 1733.3334
$ws.Range("K34").Value = 3860.25
$ws.Range("L34").Value = 5200.0002
$ws.Range("M34").Value = -3776.25
$ws.Range("N34").Value = -5368.0002

$ws.Range("H39").Value = 3799.8
$ws.Range("J39").Value = 3799.8
$ws.Range("L39").Value = 11399.4
$ws.Range("N39").Value = -11987.4

$ws.Range("H50").Value = 167445.83
$ws.Range("I50").Value = 188
$ws.Range("J50").Value = 251074.75
$ws.Range("K50").Value = 564
$ws.Range("L50").Value = 753224.25
$ws.Range("M50").Value = -83
$ws.Range("N50").Value = -754186.25

$ws.Range("H53").Value = 167445.83
$ws.Range("I53").Value = 188
$ws.Range("J53").Value = 251074.75
$ws.Range("K53").Value = 564
$ws.Range("L53").Value = 753224.25
$ws.Range("M53").Value = -83
$ws.Range("N53").Value = -754186.25

$ws.Range("H55").Value = 4098
$ws.Range("I55").Value = 349
$ws.Range("J55").Value = 7097.2
$ws.Range("K55").Value = 1047
$ws.Range("L55").Value = 21291.6
$ws.Range("M55").Value = -870
$ws.Range("N55").Value = -21645.6

$ws.Range("H56").Value = 5148.9756
$ws.Range("I56").Value = 5148.9756
$ws.Range("K56").Value = 5148.9756
$ws.Range("M56").Value = -4618.9756

$ws.Range("H63").Value = 6000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 18000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -19498

$ws.Range("H66").Value = 6000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 54000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -61488

$ws.Range("H107").Value = 334437.12
$ws.Range("J107").Value = 385785.03
$ws.Range("L107").Value = 1157355.09
$ws.Range("N107").Value = -1161195.09

$ws.Range("H122").Value = 11970.223
$ws.Range("I122").Value = 896.25
$ws.Range("J122").Value = 20829.4
$ws.Range("K122").Value = 8066.25
$ws.Range("L122").Value = 187464.6
$ws.Range("M122").Value = -5616.25
$ws.Range("N122").Value = -192364.6

$ws.Range("H132").Value = 2505.25
$ws.Range("I132").Value = 2489.5
$ws.Range("J132").Value = 2552.5
$ws.Range("K132").Value = 22405.5
$ws.Range("L132").Value = 22972.5
$ws.Range("M132").Value = -19875.5
$ws.Range("N132").Value = -28032.5

$ws.Range("H133").Value = 7218.8887
$ws.Range("I133").Value = 4996.25
$ws.Range("K133").Value = 14988.75
$ws.Range("M133").Value = -9928.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 38562
$ws.Range("J105").Value = 38562
$ws.Range("L105").Value = 38562
$ws.Range("N105").Value = -45550

$ws.Range("H140").Value = 89999
$ws.Range("J140").Value = 89999
$ws.Range("L140").Value = 89999
$ws.Range("N140").Value = -100359

$ws.Range("H141").Value = 48424.168
$ws.Range("J141").Value = 48424.168
$ws.Range("L141").Value = 48424.168
$ws.Range("N141").Value = -58784.168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3438
$ws.Range("I16").Value = 3032.3635
$ws.Range("K16").Value = 3032.3635
$ws.Range("M16").Value = -2862.3635

$ws.Range("H40").Value = 25875.75
$ws.Range("I40").Value = 16200.2
$ws.Range("K40").Value = 16200.2
$ws.Range("M40").Value = -16064.2

$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524

$ws.Range("H122").Value = 9166.666999999999
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 8750
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 26250
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -31150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 23715
$ws.Range("J105").Value = 23715
$ws.Range("L105").Value = 23715
$ws.Range("N105").Value = -30703

$ws.Range("H113").Value = 5737.5
$ws.Range("J113").Value = 6318.3335
$ws.Range("L113").Value = 18955.0005
$ws.Range("N113").Value = -23295.0005

$ws.Range("H126").Value = 72905.13
$ws.Range("I126").Value = 105597.7
$ws.Range("J126").Value = 7520
$ws.Range("K126").Value = 316793.1
$ws.Range("L126").Value = 22560
$ws.Range("M126").Value = -314323.1
$ws.Range("N126").Value = -27500
